# Fruta / hortaliza, semanal
# A new weekly data row is inserted just above the existing row 98
# (Feria Lagunitas de Puerto Montt / Pomelo), pushing all subsequent
# rows down by one (old row 143 becomes row 144). The dimension grows
# from A1:T143 to A1:T144.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 98, shifting rows 98-143 down to 99-144.
$ws.Rows.Item(98).Insert()

# Populate the newly inserted row 98 with the new weekly record.
$ws.Cells.Item(98, 1).Value2 = 4
$ws.Cells.Item(98, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(98, 3).Value2 = "Los Lagos"
$ws.Cells.Item(98, 4).Value2 = 44466
$ws.Cells.Item(98, 5).Value2 = 10
$ws.Cells.Item(98, 6).Value2 = "Fruta"
$ws.Cells.Item(98, 7).Value2 = 100102
$ws.Cells.Item(98, 8).Value2 = "Cítricos"
$ws.Cells.Item(98, 9).Value2 = 100102006
$ws.Cells.Item(98, 10).Value2 = "Pomelo"
$ws.Cells.Item(98, 11).Value2 = "Start Ruby"
$ws.Cells.Item(98, 12).Value2 = "Primera"
$ws.Cells.Item(98, 13).Value2 = 60
$ws.Cells.Item(98, 14).Value2 = 12000
$ws.Cells.Item(98, 15).Value2 = 12000
$ws.Cells.Item(98, 16).Value2 = 12000
$ws.Cells.Item(98, 17).Value2 = "`$/caja 14 kilos empedrada"
$ws.Cells.Item(98, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(98, 19).Value2 = 857
$ws.Cells.Item(98, 20).Value2 = 14
